$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it after
# touching the "Data" sheet's selection (the workbook's active tab stays
# on "Graph" in the target state).
$originalActiveSheet = $wb.ActiveSheet

# Update the "Data" sheet: 2025 Book Count 30 -> 31 (year-end count revised
# up by one book), and leave the sheet's remembered selection on D18.
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()
$ws.Range("B14").Value = 31
$ws.Range("D18").Select()

# Refresh any dependent charts/pivots now that the source data changed.
$wb.RefreshAll()

# Restore the originally active sheet/tab.
$originalActiveSheet.Activate()

$wb.Save()
